# Updated cryptos list with latest price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.242.76'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '3.128.16'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.82'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.125.04'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.17'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '3.639.70'
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("D17").Value = '64.226.14'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '3.132.46'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.89'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.93'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.55'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.68'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.22'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.40'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.75'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.36'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +7.97%  '
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.114'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.88'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.64'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("E35").Value = '  -0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.96'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").Value = '0.0₃0768'
$ws.Range("E37").Value = '  +5.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.35'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '445.04'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.20'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("D44").Value = '2.854.05'
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.22'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.43'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.00'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.34'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.01%  '
